$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2781226666666667
$ws.Range("H2").Value = 0.834368
$ws.Range("I2").Value = 0.2500965174582698
$ws.Range("J2").Value = 0.2500965174582697
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.356166333333333
$ws.Range("N2").Value = 4.068499
$ws.Range("O2").Value = 0.4438852433350371
$ws.Range("P2").Value = 0.4438852433350371
$ws.Range("Q2").Value = 0.3771805970702223
$ws.Range("R2").Value = 3.394625373632
$ws.Range("S2").Value = 0.1110141535092094
$ws.Range("T2").Value = 0.1110141535092094
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2781226666666667
$ws.Range("H3").Value = 0.834368
$ws.Range("I3").Value = 0.2500965174582698
$ws.Range("J3").Value = 0.2500965174582697
$ws.Range("M3").Value = 0.9609030000000001
$ws.Range("O3").Value = 0.3145120561487422
$ws.Range("P3").Value = 0.3145120561487422
$ws.Range("Q3").Value = 0.2672489047680001
$ws.Range("R3").Value = 2.405240142912
$ws.Range("S3").Value = 0.07865836994144024
$ws.Range("T3").Value = 0.07865836994144021
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2781226666666667
$ws.Range("H4").Value = 0.834368
$ws.Range("I4").Value = 0.2500965174582698
$ws.Range("J4").Value = 0.2500965174582697
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.525837
$ws.Range("N4").Value = 1.577511
$ws.Range("O4").Value = 0.1721111038981938
$ws.Range("P4").Value = 0.1721111038981938
$ws.Range("Q4").Value = 0.146247188672
$ws.Range("R4").Value = 1.316224698048
$ws.Range("S4").Value = 0.04304438770083672
$ws.Range("T4").Value = 0.04304438770083671
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2781226666666667
$ws.Range("H5").Value = 0.834368
$ws.Range("I5").Value = 0.2500965174582698
$ws.Range("J5").Value = 0.2500965174582697
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.212312
$ws.Range("N5").Value = 0.6369359999999999
$ws.Range("O5").Value = 0.06949159661802674
$ws.Range("P5").Value = 0.06949159661802674
$ws.Range("Q5").Value = 0.05904877960533333
$ws.Range("R5").Value = 0.531439016448
$ws.Range("S5").Value = 0.01737960630678336
$ws.Range("T5").Value = 0.01737960630678336
$ws.Range("I6").Value = 0.5435053941868915
$ws.Range("J6").Value = 0.5435053941868914
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.356166333333333
$ws.Range("N6").Value = 4.068499
$ws.Range("O6").Value = 0.4438852433350371
$ws.Range("P6").Value = 0.4438852433350371
$ws.Range("Q6").Value = 0.8196823017517777
$ws.Range("R6").Value = 7.377140715766
$ws.Range("S6").Value = 0.2412540241525536
$ws.Range("T6").Value = 0.2412540241525535
$ws.Range("I7").Value = 0.5435053941868915
$ws.Range("J7").Value = 0.5435053941868914
$ws.Range("M7").Value = 0.9609030000000001
$ws.Range("O7").Value = 0.3145120561487422
$ws.Range("P7").Value = 0.3145120561487422
$ws.Range("Q7").Value = 0.580780663434
$ws.Range("R7").Value = 5.227025970906
$ws.Range("S7").Value = 0.1709389990536519
$ws.Range("T7").Value = 0.1709389990536518
$ws.Range("I8").Value = 0.5435053941868915
$ws.Range("J8").Value = 0.5435053941868914
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.525837
$ws.Range("N8").Value = 1.577511
$ws.Range("O8").Value = 0.1721111038981938
$ws.Range("P8").Value = 0.1721111038981938
$ws.Range("Q8").Value = 0.317821842286
$ws.Range("R8").Value = 2.860396580574
$ws.Range("S8").Value = 0.09354331336812888
$ws.Range("T8").Value = 0.09354331336812886
$ws.Range("I9").Value = 0.5435053941868915
$ws.Range("J9").Value = 0.5435053941868914
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.212312
$ws.Range("N9").Value = 0.6369359999999999
$ws.Range("O9").Value = 0.06949159661802674
$ws.Range("P9").Value = 0.06949159661802674
$ws.Range("Q9").Value = 0.1283237790026666
$ws.Range("R9").Value = 1.154914011024
$ws.Range("S9").Value = 0.03776905761255708
$ws.Range("T9").Value = 0.03776905761255707
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.183755
$ws.Range("H10").Value = 0.551265
$ws.Range("I10").Value = 0.1652381883013647
$ws.Range("J10").Value = 0.1652381883013647
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.356166333333333
$ws.Range("N10").Value = 4.068499
$ws.Range("O10").Value = 0.4438852433350371
$ws.Range("P10").Value = 0.4438852433350371
$ws.Range("Q10").Value = 0.2492023445816667
$ws.Range("R10").Value = 2.242821101235
$ws.Range("S10").Value = 0.07334679342239196
$ws.Range("T10").Value = 0.07334679342239195
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.183755
$ws.Range("H11").Value = 0.551265
$ws.Range("I11").Value = 0.1652381883013647
$ws.Range("J11").Value = 0.1652381883013647
$ws.Range("M11").Value = 0.9609030000000001
$ws.Range("O11").Value = 0.3145120561487422
$ws.Range("P11").Value = 0.3145120561487422
$ws.Range("Q11").Value = 0.176570730765
$ws.Range("R11").Value = 1.589136576885
$ws.Range("S11").Value = 0.05196940235695526
$ws.Range("T11").Value = 0.05196940235695525
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.183755
$ws.Range("H12").Value = 0.551265
$ws.Range("I12").Value = 0.1652381883013647
$ws.Range("J12").Value = 0.1652381883013647
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.525837
$ws.Range("N12").Value = 1.577511
$ws.Range("O12").Value = 0.1721111038981938
$ws.Range("P12").Value = 0.1721111038981938
$ws.Range("Q12").Value = 0.096625177935
$ws.Range("R12").Value = 0.869626601415
$ws.Range("S12").Value = 0.0284393269946855
$ws.Range("T12").Value = 0.02843932699468549
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.183755
$ws.Range("H13").Value = 0.551265
$ws.Range("I13").Value = 0.1652381883013647
$ws.Range("J13").Value = 0.1652381883013647
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.212312
$ws.Range("N13").Value = 0.6369359999999999
$ws.Range("O13").Value = 0.06949159661802674
$ws.Range("P13").Value = 0.06949159661802674
$ws.Range("Q13").Value = 0.03901339155999999
$ws.Range("R13").Value = 0.35112052404
$ws.Range("S13").Value = 0.01148266552733198
$ws.Range("T13").Value = 0.01148266552733198
$ws.Range("G14").Value = 0.04577233333333334
$ws.Range("H14").Value = 0.137317
$ws.Range("I14").Value = 0.04115990005347428
$ws.Range("J14").Value = 0.04115990005347428
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.356166333333333
$ws.Range("N14").Value = 4.068499
$ws.Range("O14").Value = 0.4438852433350371
$ws.Range("P14").Value = 0.4438852433350371
$ws.Range("Q14").Value = 0.06207489746477778
$ws.Range("R14").Value = 0.5586740771830001
$ws.Range("S14").Value = 0.01827027225088224
$ws.Range("T14").Value = 0.01827027225088224
$ws.Range("G15").Value = 0.04577233333333334
$ws.Range("H15").Value = 0.137317
$ws.Range("I15").Value = 0.04115990005347428
$ws.Range("J15").Value = 0.04115990005347428
$ws.Range("M15").Value = 0.9609030000000001
$ws.Range("O15").Value = 0.3145120561487422
$ws.Range("P15").Value = 0.3145120561487422
$ws.Range("Q15").Value = 0.04398277241700001
$ws.Range("R15").Value = 0.3958449517530001
$ws.Range("S15").Value = 0.01294528479669492
$ws.Range("T15").Value = 0.01294528479669492
$ws.Range("G16").Value = 0.04577233333333334
$ws.Range("H16").Value = 0.137317
$ws.Range("I16").Value = 0.04115990005347428
$ws.Range("J16").Value = 0.04115990005347428
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.525837
$ws.Range("N16").Value = 1.577511
$ws.Range("O16").Value = 0.1721111038981938
$ws.Range("P16").Value = 0.1721111038981938
$ws.Range("Q16").Value = 0.024068786443
$ws.Range("R16").Value = 0.2166190779870001
$ws.Range("S16").Value = 0.007084075834542786
$ws.Range("T16").Value = 0.007084075834542786
$ws.Range("G17").Value = 0.04577233333333334
$ws.Range("H17").Value = 0.137317
$ws.Range("I17").Value = 0.04115990005347428
$ws.Range("J17").Value = 0.04115990005347428
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.212312
$ws.Range("N17").Value = 0.6369359999999999
$ws.Range("O17").Value = 0.06949159661802674
$ws.Range("P17").Value = 0.06949159661802674
$ws.Range("Q17").Value = 0.009718015634666666
$ws.Range("R17").Value = 0.08746214071200001
$ws.Range("S17").Value = 0.002860267171354332
$ws.Range("T17").Value = 0.002860267171354332
